$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44432
$ws.Range("M2").Value = 10
$ws.Range("R2").Value = 'Perú'
$ws.Range("D3").Value = 44424
$ws.Range("M3").Value = 15
$ws.Range("N3").Value = 35000
$ws.Range("O3").Value = 35000
$ws.Range("P3").Value = 35000
$ws.Range("S3").Value = 1944
$ws.Range("D4").Value = 44369
$ws.Range("M4").Value = 5
$ws.Range("D5").Value = 44418
$ws.Range("M5").Value = 30
$ws.Range("N5").Value = 35000
$ws.Range("O5").Value = 35000
$ws.Range("P5").Value = 35000
$ws.Range("Q5").Value = '$/caja 18 kilos'
$ws.Range("S5").Value = 1944
$ws.Range("T5").Value = 18
$ws.Range("D6").Value = 44377
$ws.Range("M6").Value = 30
$ws.Range("N6").Value = 40000
$ws.Range("O6").Value = 40000
$ws.Range("P6").Value = 40000
$ws.Range("R6").Value = 'Región de Arica y Parinacota'
$ws.Range("S6").Value = 2222
$ws.Range("D7").Value = 44357
$ws.Range("M7").Value = 10
$ws.Range("R7").Value = 'Perú'
$ws.Range("D8").Value = 44363
$ws.Range("M8").Value = 144
$ws.Range("N8").Value = 1700
$ws.Range("O8").Value = 1700
$ws.Range("P8").Value = 1700
$ws.Range("Q8").Value = '$/kilo'
$ws.Range("R8").Value = 'Región de Arica y Parinacota'
$ws.Range("S8").Value = 1700
$ws.Range("T8").Value = 1
$ws.Range("D9").Value = 44379
$ws.Range("N9").Value = 30000
$ws.Range("O9").Value = 30000
$ws.Range("P9").Value = 30000
$ws.Range("R9").Value = 'Región de Arica y Parinacota'
$ws.Range("S9").Value = 1667
$ws.Range("D10").Value = 44448
$ws.Range("M10").Value = 50
$ws.Range("N10").Value = 38000
$ws.Range("O10").Value = 38000
$ws.Range("P10").Value = 38000
$ws.Range("S10").Value = 2111
$ws.Range("D11").Value = 44294
$ws.Range("D12").Value = 44435
$ws.Range("M12").Value = 10
$ws.Range("N12").Value = 35000
$ws.Range("O12").Value = 35000
$ws.Range("P12").Value = 35000
$ws.Range("Q12").Value = '$/caja 18 kilos'
$ws.Range("R12").Value = 'Perú'
$ws.Range("S12").Value = 1944
$ws.Range("T12").Value = 18
$ws.Range("D13").Value = 44435
$ws.Range("M13").Value = 105
$ws.Range("D14").Value = 44405
$ws.Range("M14").Value = 10
$ws.Range("D15").Value = 44434
$ws.Range("M15").Value = 40
$ws.Range("D16").Value = 44433
$ws.Range("M16").Value = 15
$ws.Range("N16").Value = 35000
$ws.Range("O16").Value = 35000
$ws.Range("P16").Value = 35000
$ws.Range("R16").Value = 'Región de Arica y Parinacota'
$ws.Range("S16").Value = 1944
$ws.Range("D17").Value = 44364
$ws.Range("M17").Value = 90
$ws.Range("N17").Value = 1700
$ws.Range("O17").Value = 1700
$ws.Range("P17").Value = 1700
$ws.Range("Q17").Value = '$/kilo'
$ws.Range("S17").Value = 1700
$ws.Range("T17").Value = 1
$ws.Range("D18").Value = 44279
$ws.Range("M18").Value = 30
$ws.Range("N18").Value = 35000
$ws.Range("O18").Value = 36000
$ws.Range("P18").Value = 35667
$ws.Range("S18").Value = 1982
$ws.Range("D19").Value = 44431
$ws.Range("M19").Value = 30
$ws.Range("N19").Value = 35000
$ws.Range("O19").Value = 35000
$ws.Range("P19").Value = 35000
$ws.Range("S19").Value = 1944
$ws.Range("D20").Value = 44442
$ws.Range("R20").Value = 'Perú'
$ws.Range("D21").Value = 44264
$ws.Range("N21").Value = 40000
$ws.Range("O21").Value = 40000
$ws.Range("P21").Value = 40000
$ws.Range("S21").Value = 2222
$ws.Range("D22").Value = 44438
$ws.Range("M22").Value = 25
$ws.Range("D23").Value = 44449
$ws.Range("N23").Value = 38000
$ws.Range("O23").Value = 38000
$ws.Range("P23").Value = 38000
$ws.Range("S23").Value = 2111
$ws.Range("D24").Value = 44392
$ws.Range("M24").Value = 20
$ws.Range("O24").Value = 35000
$ws.Range("P24").Value = 35000
$ws.Range("S24").Value = 1944
